$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("L4").Value = "ICD10 GM diagnosis code mapped A = 1, B = 2, C = 3, D = 4, e.g.: A01.9 = 101.9, C50.1 = 350.1 or D41.9 = 441.9"
$ws.Range("L5").Value = "ICD10 GM diagnosis code grouped to parent code, e.g. A01.1 and A01.9 both belong to group 101 (remove decimal from icd10_mapped)"
$ws.Range("L6").Value = "ICD10 GM diagnosis code grouped to entity groups from 0-23 according to LGL Report Cancer in Bavaria 2019, mapping see github.com/bzkf/onco-analytics-on-fhir/src/obds_fhir_to_opal/utils_onco_analytics.py"
$ws.Range("L7").Value = "date of diagnosis"
$ws.Range("L9").Value = "Month of diagnosis"
$ws.Range("L10").Value = "Day of diagnosis"
$ws.Range("L12").Value = "Gender mapped: 0 = None, 1 = female, 2 = male, 3 = other/diverse"
